$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of table row index -> value to insert into the (empty) second cell of that row.
# Row 24 = "Ratio" (Retention ratio), rows 44-46 = Answer Recall Lenient / Strict / Average.
$values = @{
    24 = "0.3333";
    44 = "0.2857";
    45 = "0.1428";
    46 = "0.2142"
}

foreach ($rowIndex in $values.Keys) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(2)
    $cell.Range.Text = $values[$rowIndex]
    $cell.Range.Font.Bold = $true
    $cell.Range.Font.Size = 12
    $cell.Range.Font.SizeBi = 12
}
